$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new value in cell A2 (references the new shared string "US1")
$ws.Range("A2").Value = "US1"

# Update the selection to reflect the author's last active cell (C10)
$ws.Range("C10").Select()
